$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.39010000000002
$ws.Range("E5").Value = 12.97549999999999
$ws.Range("E9").Value = 13.38710000000001
$ws.Range("E11").Value = 13.37119999999999
$ws.Range("A21").Value = -21.1437
$ws.Range("E21").Value = 12.74800000000001
$ws.Range("A23").Value = -21.34420000000002
$ws.Range("A25").Value = -22.44860000000003
